$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.59528094111439
$ws.Range("C2").Value = 5.446275288111754
$ws.Range("E2").Value = 16.42572648880487
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 32.59227319600827
$ws.Range("H2").Value = 15.06553459527227
$ws.Range("K2").Value = 9.82829173266871
$ws.Range("N2").Value = 18.54925090105182
$ws.Range("B3").Value = 10.2867226849991
$ws.Range("C3").Value = 5.20473591566587
$ws.Range("E3").Value = 15.49831418698885
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 32.51606518544829
$ws.Range("H3").Value = 15.10828100784827
$ws.Range("K3").Value = 9.61393358514103
$ws.Range("N3").Value = 18.62029197098168
$ws.Range("B4").Value = 10.09535070360622
$ws.Range("C4").Value = 5.049469065128872
$ws.Range("E4").Value = 14.90462964268383
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 32.48158791785345
$ws.Range("H4").Value = 15.13761027031671
$ws.Range("K4").Value = 9.482419914307274
$ws.Range("N4").Value = 18.66584579338654
$ws.Range("B5").Value = 10.0170156362008
$ws.Range("C5").Value = 4.984490625622399
$ws.Range("E5").Value = 14.65687243619023
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 32.47063311013577
$ws.Range("H5").Value = 15.15033490889851
$ws.Range("K5").Value = 9.428934765716122
$ws.Range("N5").Value = 18.68489728848311
$ws.Range("B6").Value = 10.00399106641571
$ws.Range("C6").Value = 4.973599409822829
$ws.Range("E6").Value = 14.61538965735225
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 32.46900086772143
$ws.Range("H6").Value = 15.15249442291738
$ws.Range("K6").Value = 9.420062531296672
$ws.Range("N6").Value = 18.68809029057998
$ws.Range("B7").Value = 10.09429548398489
$ws.Range("C7").Value = 5.048599586487948
$ws.Range("E7").Value = 14.90131149321191
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 32.48142765169382
$ws.Range("H7").Value = 15.13777875387051
$ws.Range("K7").Value = 9.481698049230626
$ws.Range("N7").Value = 18.66610075085097
$ws.Range("B8").Value = 10.48936474862728
$ws.Range("C8").Value = 5.364461671333793
$ws.Range("E8").Value = 16.11113517942834
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 32.56344101360031
$ws.Range("H8").Value = 15.07963225761446
$ws.Range("K8").Value = 9.754408107862897
$ws.Range("N8").Value = 18.57334532802086
$ws.Range("B9").Value = 11.243479979921
$ws.Range("C9").Value = 5.927027970871957
$ws.Range("E9").Value = 18.34456681579281
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 32.82185842347321
$ws.Range("H9").Value = 14.99017634773178
$ws.Range("K9").Value = 10.28652981058575
$ws.Range("N9").Value = 18.4067282550447
$ws.Range("B10").Value = 11.77832176481347
$ws.Range("C10").Value = 6.303956343374518
$ws.Range("E10").Value = 19.97768788970688
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 33.07076166922681
$ws.Range("H10").Value = 14.93957713440493
$ws.Range("K10").Value = 10.67140587729016
$ws.Range("N10").Value = 18.29352453032176
$ws.Range("B11").Value = 12.01624075562886
$ws.Range("C11").Value = 6.467261492619055
$ws.Range("E11").Value = 20.67881211414785
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 33.19662909873117
$ws.Range("H11").Value = 14.91987227630466
$ws.Range("K11").Value = 10.8443072271304
$ws.Range("N11").Value = 18.24400303805143
$ws.Range("B12").Value = 12.10547067470369
$ws.Range("C12").Value = 6.527910091434188
$ws.Range("E12").Value = 20.938337035268
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 33.2460871488336
$ws.Range("H12").Value = 14.9128893017747
$ws.Range("K12").Value = 10.90940003272977
$ws.Range("N12").Value = 18.2255329800916
$ws.Range("B13").Value = 12.08629323027129
$ws.Range("C13").Value = 6.51490157271466
$ws.Range("E13").Value = 20.88270892996426
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 33.23535605579701
$ws.Range("H13").Value = 14.91437187657428
$ws.Range("K13").Value = 10.8953991144226
$ws.Range("N13").Value = 18.22949828865894
$ws.Range("B14").Value = 12.02359957585927
$ws.Range("C14").Value = 6.472275059599284
$ws.Range("E14").Value = 20.70028295251247
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 33.20066221729594
$ws.Range("H14").Value = 14.91928817278544
$ws.Range("K14").Value = 10.84967046222993
$ws.Range("N14").Value = 18.24247783845579
$ws.Range("B15").Value = 11.98508273231952
$ws.Range("C15").Value = 6.446009445722336
$ws.Range("E15").Value = 20.58776457455382
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 33.1796442265484
$ws.Range("H15").Value = 14.92236197416813
$ws.Range("K15").Value = 10.82160877094863
$ws.Range("N15").Value = 18.25046495600283
$ws.Range("B16").Value = 11.76265691611473
$ws.Range("C16").Value = 6.293118127950457
$ws.Range("E16").Value = 19.93102973485188
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 33.06278842791031
$ws.Range("H16").Value = 14.94093176807136
$ws.Range("K16").Value = 10.66005660597531
$ws.Range("N16").Value = 18.29680050653977
$ws.Range("B17").Value = 11.62475972623316
$ws.Range("C17").Value = 6.197220046620273
$ws.Range("E17").Value = 19.51746472535995
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 32.99432365787069
$ws.Range("H17").Value = 14.95317402081125
$ws.Range("K17").Value = 10.56034155189111
$ws.Range("N17").Value = 18.32573081871879
$ws.Range("B18").Value = 11.54494306740361
$ws.Range("C18").Value = 6.141294642570876
$ws.Range("E18").Value = 19.27565751831151
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 32.9561359101499
$ws.Range("H18").Value = 14.96052713391502
$ws.Range("K18").Value = 10.50278659796059
$ws.Range("N18").Value = 18.34255679126056
$ws.Range("B19").Value = 11.51783523657579
$ws.Range("C19").Value = 6.122228069290953
$ws.Range("E19").Value = 19.19310878275819
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 32.9434114444341
$ws.Range("H19").Value = 14.96307023004627
$ws.Range("K19").Value = 10.48326700550979
$ws.Range("N19").Value = 18.34828576916837
$ws.Range("B20").Value = 11.63949174166996
$ws.Range("C20").Value = 6.207508121112555
$ws.Range("E20").Value = 19.56189649573813
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 33.00148869886643
$ws.Range("H20").Value = 14.9518385351402
$ws.Range("K20").Value = 10.57097774709443
$ws.Range("N20").Value = 18.32263189521039
$ws.Range("B21").Value = 12.04203837281727
$ws.Range("C21").Value = 6.484827963062171
$ws.Range("E21").Value = 20.75402779168701
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 33.21080413073251
$ws.Range("H21").Value = 14.91783112326509
$ws.Range("K21").Value = 10.86311293765009
$ws.Range("N21").Value = 18.23865777002333
$ws.Range("B22").Value = 12.30004897673694
$ws.Range("C22").Value = 6.659121438335227
$ws.Range("E22").Value = 21.49833965841757
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 33.35805184831899
$ws.Range("H22").Value = 14.8983973510091
$ws.Range("K22").Value = 11.05179016122435
$ws.Range("N22").Value = 18.18542271599629
$ws.Range("B23").Value = 12.16283611394691
$ws.Range("C23").Value = 6.566738811266752
$ws.Range("E23").Value = 21.10426083909906
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 33.27851556009056
$ws.Range("H23").Value = 14.9085132705864
$ws.Range("K23").Value = 10.95131662669635
$ws.Range("N23").Value = 18.21368503887624
$ws.Range("B24").Value = 11.63283306536188
$ws.Range("C24").Value = 6.202859349263517
$ws.Range("E24").Value = 19.54182148932415
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 32.99824572633592
$ws.Range("H24").Value = 14.9524413274609
$ws.Range("K24").Value = 10.5661698264306
$ws.Range("N24").Value = 18.32403231543167
$ws.Range("B25").Value = 11.04239903897623
$ws.Range("C25").Value = 5.781123642674929
$ws.Range("E25").Value = 17.71725841758899
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 32.7415312460634
$ws.Range("H25").Value = 15.01173089412739
$ws.Range("K25").Value = 10.14333475365413
$ws.Range("N25").Value = 18.45017785440143
